# Memperbarui coding utama dan membuat struk
# Append three new response rows (12-14) to the form-submissions sheet,
# matching the data captured after rows 1-11.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12
$ws.Range("A12").Value = 12
$ws.Range("B12").Value = "btara"
$ws.Range("C12").Value = "'123"
$ws.Range("D12").Value = "13/06/2022 10:40:05"

# Row 13
$ws.Range("A13").Value = 13
$ws.Range("B13").Value = "btara"
$ws.Range("C13").Value = "'123"
$ws.Range("D13").Value = "14/06/2022 12:05:49"

# Row 14
$ws.Range("A14").Value = 14
$ws.Range("B14").Value = "aaa"
$ws.Range("C14").Value = "aaa"
$ws.Range("D14").Value = "14/06/2022 19:16:43"
